$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing data (shifts LastName/Company/Expected right)
$ws.Range("A1:B1").EntireColumn.Insert()

# New data row values
$ws.Range("A2").Value = "hari.radhakrishnan@qeagle.com"
$ws.Range("B2").Value = "Testleaf`$321"

# New header row values
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Add hyperlink on A2
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:hari.radhakrishnan@qeagle.com")

[void]$ws.Range("B1").Select()
